# Auto-generated: apply updated market-price / profit figures to each sheet
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1297.2858
$ws.Range("I41").Value = 1530.091
$ws.Range("K41").Value = 1530.091
$ws.Range("M41").Value = -1090.091
$ws.Range("H88").Value = 4049
$ws.Range("J88").Value = 4049
$ws.Range("L88").Value = 4049
$ws.Range("N88").Value = -4861
$ws.Range("H91").Value = 4049
$ws.Range("J91").Value = 4049
$ws.Range("L91").Value = 4049
$ws.Range("N91").Value = -6857
$ws.Range("H103").Value = 446.6
$ws.Range("I103").Value = 465.5
$ws.Range("K103").Value = 1396.5
$ws.Range("M103").Value = -810.5
$ws.Range("H132").Value = 4318
$ws.Range("I132").Value = 4435.3774
$ws.Range("J132").Value = 2762.75
$ws.Range("K132").Value = 13306.1322
$ws.Range("L132").Value = 8288.25
$ws.Range("M132").Value = -10776.1322
$ws.Range("N132").Value = -13348.25
$ws.Range("H137").Value = 1923.6666
$ws.Range("I137").Value = 1882.3077
$ws.Range("K137").Value = 5646.9231
$ws.Range("M137").Value = -3096.9231
$ws.Range("H138").Value = 6607.82
$ws.Range("I138").Value = 4953.4
$ws.Range("J138").Value = 6899.7764
$ws.Range("K138").Value = 14860.2
$ws.Range("L138").Value = 20699.3292
$ws.Range("M138").Value = -9720.199999999999
$ws.Range("N138").Value = -30979.3292
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5413.8423
$ws.Range("I61").Value = 3222.8572
$ws.Range("J61").Value = 11548.6
$ws.Range("K61").Value = 3222.8572
$ws.Range("L61").Value = 11548.6
$ws.Range("M61").Value = -3010.8572
$ws.Range("N61").Value = -11972.6
$ws.Range("H110").Value = 1485.3
$ws.Range("J110").Value = 3000
$ws.Range("L110").Value = 3000
$ws.Range("N110").Value = -7090
$ws.Range("H132").Value = 3613.7693
$ws.Range("I132").Value = 2845.5264
$ws.Range("K132").Value = 8536.5792
$ws.Range("M132").Value = -6006.5792
$ws.Range("H136").Value = 5413.8423
$ws.Range("I136").Value = 3222.8572
$ws.Range("J136").Value = 11548.6
$ws.Range("K136").Value = 9668.571599999999
$ws.Range("L136").Value = 34645.8
$ws.Range("M136").Value = -7118.571599999999
$ws.Range("N136").Value = -39745.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 21187.2
$ws.Range("I94").Value = 1484
$ws.Range("K94").Value = 1484
$ws.Range("M94").Value = -1033
$ws.Range("H105").Value = 2155.3333
$ws.Range("I105").Value = 1859.2727
$ws.Range("K105").Value = 1859.2727
$ws.Range("M105").Value = -112.2727
$ws.Range("H134").Value = 8479.091
$ws.Range("I134").Value = 8284.846
$ws.Range("K134").Value = 24854.538
$ws.Range("M134").Value = -22319.538
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 929
$ws.Range("I16").Value = 786.25
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 786.25
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -499.25
$ws.Range("N16").Value = -2074
$ws.Range("H31").Value = 38950.715
$ws.Range("I31").Value = 3184.4783
$ws.Range("J31").Value = 203475.4
$ws.Range("K31").Value = 3184.4783
$ws.Range("L31").Value = 203475.4
$ws.Range("M31").Value = -2889.4783
$ws.Range("N31").Value = -204065.4
$ws.Range("H34").Value = 38950.715
$ws.Range("I34").Value = 3184.4783
$ws.Range("J34").Value = 203475.4
$ws.Range("K34").Value = 3184.4783
$ws.Range("L34").Value = 203475.4
$ws.Range("M34").Value = -2982.4783
$ws.Range("N34").Value = -203879.4
$ws.Range("H58").Value = 4190.5625
$ws.Range("I58").Value = 3817.7144
$ws.Range("K58").Value = 3817.7144
$ws.Range("M58").Value = -3614.7144
$ws.Range("H99").Value = 2446.4092
$ws.Range("I99").Value = 2373.2856
$ws.Range("J99").Value = 2574.375
$ws.Range("K99").Value = 2373.2856
$ws.Range("L99").Value = 2574.375
$ws.Range("M99").Value = -875.2856000000002
$ws.Range("N99").Value = -5570.375
$ws.Range("H113").Value = 929
$ws.Range("I113").Value = 786.25
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 786.25
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1383.75
$ws.Range("N113").Value = -5840
$ws.Range("H126").Value = 2446.4092
$ws.Range("I126").Value = 2373.2856
$ws.Range("J126").Value = 2574.375
$ws.Range("K126").Value = 7119.8568
$ws.Range("L126").Value = 7723.125
$ws.Range("M126").Value = -4649.8568
$ws.Range("N126").Value = -12663.125
$ws.Range("H136").Value = 4190.5625
$ws.Range("I136").Value = 3817.7144
$ws.Range("K136").Value = 11453.1432
$ws.Range("M136").Value = -8903.143199999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7017.5
$ws.Range("I56").Value = 7017.5
$ws.Range("K56").Value = 7017.5
$ws.Range("M56").Value = -6487.5
$ws.Range("I63").Value = 1000
$ws.Range("K63").Value = 3000
$ws.Range("M63").Value = -2251
$ws.Range("I66").Value = 1000
$ws.Range("K66").Value = 9000
$ws.Range("M66").Value = -5256
$ws.Range("H75").Value = 2498.5
$ws.Range("J75").Value = 2498.5
$ws.Range("L75").Value = 7495.5
$ws.Range("N75").Value = -9491.5
$ws.Range("H78").Value = 2498.5
$ws.Range("J78").Value = 2498.5
$ws.Range("L78").Value = 22486.5
$ws.Range("N78").Value = -32470.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 10000.5
$ws.Range("I19").Value = 10000.5
$ws.Range("K19").Value = 10000.5
$ws.Range("M19").Value = -9712.5
$ws.Range("H23").Value = 15005.5
$ws.Range("I23").Value = 12
$ws.Range("J23").Value = 29999
$ws.Range("K23").Value = 12
$ws.Range("L23").Value = 29999
$ws.Range("M23").Value = 211
$ws.Range("N23").Value = -30445
$ws.Range("H102").Value = 13431.174
$ws.Range("I102").Value = 15200.85
$ws.Range("K102").Value = 15200.85
$ws.Range("M102").Value = -13578.85
$ws.Range("H132").Value = 10165
$ws.Range("I132").Value = 4823.8667
$ws.Range("J132").Value = 21610.285
$ws.Range("K132").Value = 14471.6001
$ws.Range("L132").Value = 64830.855
$ws.Range("M132").Value = -11941.6001
$ws.Range("N132").Value = -69890.855
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3710.2666
$ws.Range("I132").Value = 3458.875
$ws.Range("J132").Value = 4715.8335
$ws.Range("K132").Value = 10376.625
$ws.Range("L132").Value = 14147.5005
$ws.Range("M132").Value = -7846.625
$ws.Range("N132").Value = -19207.5005
$ws.Range("H136").Value = 3921
$ws.Range("I136").Value = 4082.611
$ws.Range("J136").Value = 3656.5454
$ws.Range("K136").Value = 12247.833
$ws.Range("L136").Value = 10969.6362
$ws.Range("M136").Value = -9697.832999999999
$ws.Range("N136").Value = -16069.6362
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 11009.5
$ws.Range("I136").Value = 10483.667
$ws.Range("J136").Value = 13939.143
$ws.Range("K136").Value = 31451.001
$ws.Range("L136").Value = 41817.429
$ws.Range("M136").Value = -28901.001
$ws.Range("N136").Value = -46917.429
